# fix: units in data file
#
# - rename the sheet "Лист1" -> "Costs"
# - column C (scooter_lifetime_km) was using the currency number format;
#   switch it to a plain 2-decimal numeric format for the data rows
# - leave the selection where the author left it (D7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Costs"

# Fix units: scooter_lifetime_km (column C) should be a plain number, not currency
$ws.Range("C2:C4").NumberFormat = "0.00"

# Restore the cursor/selection position recorded in the saved file
$ws.Range("D7").Select()
